$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 4) to the job posting table.
# Column A reuses the same Job_Description text already stored in A3
# (the "Landing Gear System Engineer" posting), columns B/C hold the
# min/max years of experience for this additional entry.
$ws.Range("A4").Value2 = $ws.Range("A3").Value2
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 15

# Writing the long wrapped text auto-expands the row height; re-run AutoFit
# so row 4 ends up at the sheet's default height, just like the other rows.
$ws.Rows.Item(4).AutoFit()
